$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 622, pushing the existing row 622 (and everything
# below it) down by one. This matches the diff: a new data point for
# 2026/01/14 is inserted, and the rest of the daily-ranking log shifts down.
$ws.Rows.Item(622).Insert()

# Populate the newly inserted row with the new data point.
# Column A holds a text date ("2026/01/14", not a real Excel date value) and
# column B a single weekday kanji character, matching every other row in the
# sheet. Pre-formatting as text ("@") keeps Excel from auto-converting the
# "2026/01/14" string into a date serial number; ClearFormats() afterwards
# removes that temporary number-format override so the cell ends up with the
# same default (unstyled) look as its neighbors.
$ws.Cells.Item(622, 1).NumberFormat = "@"
$ws.Cells.Item(622, 1).Value = "2026/01/14"
$ws.Cells.Item(622, 2).Value = "水"
$ws.Cells.Item(622, 3).Value = 14
$ws.Cells.Item(622, 4).Value = 25
$ws.Range("A622:B622").ClearFormats()
